# test-upload-excel.xlsx : add a new student row ("Ahmad") to the table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- append the new data row (row 5) ------------------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Ahmad"
$ws.Range("C5").Value = 1710998
$ws.Range("D5").Value = "Computer Science"
$ws.Range("E5").Value = "Software Engineering"
$ws.Range("F5").Value = 2
$ws.Range("F5").NumberFormat = "0.00_);[Red]\(0.00\)"

# --- bump the sheet's recorded outline depth (row 3->4, keep col at 5) --
# touch a scratch row/col far outside the used range to advance the
# sheet's tracked outline levels, then remove the scratch row/col again
# so only the sheetFormatPr bookkeeping is left behind.
$ws.Rows.Item(100).OutlineLevel = 4
$ws.Columns.Item(50).OutlineLevel = 5
$ws.Rows.Item(100).Delete()
$ws.Columns("45:55").Delete()

# --- move the active selection on to the next empty row -----------------
$ws.Range("F6").Select()
